$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix 1: Skitty's moves (row 47) were wrongly sharing Meowsy's "Covet" text;
#     give Skitty the correct "Tackle, Growl" moves, and expand Meowsy's (row 46)
#     moves to "Covet, Growl".
$ws.Range("E47").Value = "Tackle, Growl"
$ws.Range("E46").Value = "Covet, Growl"

# --- Fix 2: append a new trainer-class block after TRAINER_MARY's Cubone entry
#     (row 112), leaving row 117 blank and moving the sheet-ending "END" marker
#     down to row 118.
$ws.Rows("113:116").Insert() | Out-Null

$ws.Range("A113").Value = ".trainerClass "
$ws.Range("B113").Value = " TRAINER_CLASS_BATTLE_GIRL,"
$ws.Range("A114").Value = ".encounterMusic_gender "
$ws.Range("B114").Value = " F_TRAINER_FEMALE | TRAINER_ENCOUNTER_MUSIC_INTENSE,"
$ws.Range("A115").Value = ".trainerPic "
$ws.Range("B115").Value = " TRAINER_PIC_BATTLE_GIRL,"
$ws.Range("A116").Value = ".items"
$ws.Range("B116").Value = "{},"

# Match the author's final cursor position.
$ws.Range("A116").Select() | Out-Null
